# Apply the Sprint 2 retrospective edits:
#  1. Slide 1: append text to the "Walter r - " run.
#  2. Slide 3: expand the "Walter -" bullet with several extra runs.

$p = $ppt.ActivePresentation
$enDash = [char]0x2013

# ---------------------------------------------------------------------
# Slide 1, Shape 1: "...Walter – " -> "...Walter – We were able to lay
# the ground work in the application to emit messages for the server."
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(1)
$tr1 = $sh1.TextFrame.TextRange

$full1 = $tr1.Text
$idx1 = $full1.IndexOf("Walter")
$run1 = $tr1.Characters($idx1 + 6, 4)   # the "r - " run right after "Walte"
$run1.Text = "r " + $enDash + " We were able to lay the ground work in the application to emit messages for the server."

# ---------------------------------------------------------------------
# Slide 3, Shape 1: "Walter -" -> "Walter - Updating team work outside
# of the ZenHub Board (ie. updating everyone) "
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(1)
$tr3 = $sh3.TextFrame.TextRange

$full3 = $tr3.Text
$idxW = $full3.IndexOf("Walter")
$runW = $tr3.Characters($idxW + 1, 8)   # the whole "Walter -" run
$runW.Text = "Walter " + $enDash + " Updating team work outside of the ZenHub Board (ie. updating everyone) "

# Re-split the new text into separate runs (matching the authored
# formatting boundaries) by re-assigning each sub-range to its own
# (unchanged) text - this preserves formatting while creating new runs.
$t = $tr3.Text

$zIdx = $t.IndexOf("ZenHub")
$rZenHub = $tr3.Characters($zIdx + 1, 6)
$rZenHub.Text = "ZenHub"

$t = $tr3.Text
$boardIdx = $t.IndexOf(" Board (", $zIdx)
$rBoard = $tr3.Characters($boardIdx + 1, 8)
$rBoard.Text = " Board ("

$t = $tr3.Text
$ieIdx = $t.IndexOf("ie", $boardIdx)
$rIe = $tr3.Characters($ieIdx + 1, 2)
$rIe.Text = "ie"

$t = $tr3.Text
$dotIdx = $t.IndexOf(". updating ", $ieIdx)
$rDot = $tr3.Characters($dotIdx + 1, 11)
$rDot.Text = ". updating "

$t = $tr3.Text
$everyoneIdx = $t.IndexOf("everyone) ", $dotIdx)
$rEveryone = $tr3.Characters($everyoneIdx + 1, 10)
$rEveryone.Text = "everyone) "
